$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells store prices and percentage changes as literal text
# (inline strings), e.g. "245.28" or "-0.70%". Excel will silently convert
# such numeric-looking strings to real numbers (and percentages to their
# decimal fraction) unless the cell is explicitly formatted as Text ("@")
# beforehand. Apply the Text format to each target cell individually, then
# write the updated value so the literal text is preserved exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.70%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.27%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.255"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.74%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05701"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.51%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.630"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.81%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.203"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.38%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8505"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.74%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9060"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.34%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1370"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.33%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07074"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.14%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.85%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09222"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.69%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001522"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.40%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005987"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.20%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005927"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.21%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.07%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.62%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.47%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03309"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.38%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1283"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.54%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.525"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04086"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.19%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.04%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001223"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.16%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.13%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.79%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03818"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.75%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1066"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.51%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.42%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002490"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.59%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01040"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "22.64%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005274"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.45%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.01%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "62.34%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002270"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.27%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
